$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: "MTF oggettivo" label and value
$ws.Range("A23").Value = "MTF oggettivo"
$ws.Range("B23").Value = 0.32

# New formulas in column G (rows 17-19)
$ws.Range("G19").Formula = "=B23*B20"
$ws.Range("G18").Formula = "=G19/B16"

# G18 picked up an auto number-format from its precedent (B16); restore the
# plain style used by the rest of the sheet by copying formats from E18
# before any other formula references G18 and inherits the same quirk.
$ws.Range("E18").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G17").Formula = "=G18/100"

# Selection ends up on G18 after these edits
[void]$ws.Range("G18").Select()
